# Update "想去人数" (F column) counts and the Cover image URL (I column)
# for matching event rows on both the "展览" and "全部类型" sheets.
#
# The two sheets list (mostly) the same events but "全部类型" has a few
# extra rows interspersed, so the row numbers between the two sheets are
# not related by a constant offset - each row is mapped explicitly below.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Each entry: 展览-row, 全部类型-row, new F value (想去人数)
$rowMap = @(
    @(7,  10, 2835),
    @(8,  11, 1681),
    @(9,  12, 1836),
    @(10, 13, 316),
    @(11, 14, 284),
    @(12, 15, 751),
    @(13, 17, 901),
    @(14, 18, 170),
    @(16, 20, 1117),
    @(20, 23, 6686),
    @(22, 25, 1553),
    @(24, 28, 180),
    @(26, 30, 310),
    @(29, 33, 1103),
    @(32, 36, 93),
    @(34, 38, 782),
    @(35, 39, 1498),
    @(36, 40, 160),
    @(37, 41, 143),
    @(38, 42, 221),
    @(39, 43, 23),
    @(41, 45, 210)
)

foreach ($entry in $rowMap) {
    $row1 = $entry[0]
    $row4 = $entry[1]
    $newVal = $entry[2]
    $ws1.Range("F$row1").Value = $newVal
    $ws4.Range("F$row4").Value = $newVal
}

# Update the Cover (I column) URL for the COMIC TIME row
# (row 7 on 展览, row 10 on 全部类型)
$newCover = "//i1.hdslb.com/bfs/openplatform/202402/9D9tHKOL1709091756797.jpeg"
$ws1.Range("I7").Value = $newCover
$ws4.Range("I10").Value = $newCover
